$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Session" to "Neurology"
$ws.Name = "Neurology"

# Append a new log row (row 42) - keep all values as text, matching the
# existing rows' storage (numberStoredAsText style data).
$row = 42

# Column A ("212024") looks like a pure number, so it would otherwise be
# auto-converted to a numeric value; force text storage for it only, then
# clear the temporary formatting so the cell keeps the sheet's default style.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "212024"
$ws.Range("A$row").ClearFormats()

# The remaining columns are not auto-detected as numbers/dates by the
# engine, so a plain value assignment already keeps them as text without
# touching their cell style.
$ws.Range("B$row").Value = "Neurology"
$ws.Range("C$row").Value = "16/12/2025"
$ws.Range("D$row").Value = "11:12:55"
$ws.Range("E$row").Value = "Scan"
$ws.Range("F$row").Value = "emp17.farah.a.youssef@gmail.com"
